# Update "Pais" dashboard sheet with new country case data (countries &
# provincias Spain update). The underlying web export re-sorts rows by
# total cases, which shifts some neighboring rows; values below reflect
# the resulting cell-by-cell content after that re-sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 09:52"

# Rusia (row 13) - updated case numbers
$ws.Range("B13").Value = 52763
$ws.Range("C13").Value = 5642
$ws.Range("D13").Value = 3873
$ws.Range("E13").Value = 48434
$ws.Range("F13").Value = 700
$ws.Range("G13").Value = 51
$ws.Range("H13").Value = 456

# Singapur moves up (new data), Rumania & Mexico shift down one row
$ws.Range("A33").Value = "Singapur"
$ws.Range("B33").Value = 9125
$ws.Range("C33").Value = 1111
$ws.Range("D33").Value = 801
$ws.Range("E33").Value = 8313
$ws.Range("F33").Value = 23
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 11

$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 8936
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 2017
$ws.Range("E34").Value = 6437
$ws.Range("F34").Value = 261
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 482

$ws.Range("A35").Value = "Mexico"
$ws.Range("B35").Value = 8772
$ws.Range("C35").Value = 511
$ws.Range("D35").Value = 2627
$ws.Range("E35").Value = 5433
$ws.Range("F35").Value = 378
$ws.Range("G35").Value = 26
$ws.Range("H35").Value = 712

# Letonia (row 90) - updated case numbers
$ws.Range("B90").Value = 748
$ws.Range("C90").Value = 9
$ws.Range("D90").Value = 88
$ws.Range("E90").Value = 655
$ws.Range("F90").Value = 3
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 5

# Sri Lanka (row 115) - updated active/recovered numbers
$ws.Range("D115").Value = 100
$ws.Range("E115").Value = 202

# Maldivas moves up (new data), Liechtenstein / Guinea Ecuatorial / Barbados shift down one row
$ws.Range("A148").Value = "Maldivas"
$ws.Range("B148").Value = 82
$ws.Range("C148").Value = 13
$ws.Range("D148").Value = 16
$ws.Range("E148").Value = 66
$ws.Range("F148").Value = 1
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 0

$ws.Range("A149").Value = "Liechtenstein"
$ws.Range("B149").Value = 81
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 55
$ws.Range("E149").Value = 25
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 1

$ws.Range("A150").Value = "Guinea Ecuatorial"
$ws.Range("B150").Value = 79
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 7
$ws.Range("E150").Value = 72
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 0

$ws.Range("A151").Value = "Barbados"
$ws.Range("B151").Value = 75
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 23
$ws.Range("E151").Value = 47
$ws.Range("F151").Value = 4
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 5
